$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore C10 to the value from the referenced revision (was 18, now 1)
$ws.Range("C10").Value = 1
